$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H28").Value = 671.087
$ws_ALC.Range("I28").Value = 778.1579
$ws_ALC.Range("J28").Value = 162.5
$ws_ALC.Range("K28").Value = 778.1579
$ws_ALC.Range("L28").Value = 162.5
$ws_ALC.Range("M28").Value = -293.1579
$ws_ALC.Range("N28").Value = -1132.5

$ws_ALC.Range("H62").Value = 3186.6667
$ws_ALC.Range("I62").Value = 3271.4285
$ws_ALC.Range("J62").Value = 2000
$ws_ALC.Range("K62").Value = 3271.4285
$ws_ALC.Range("L62").Value = 2000
$ws_ALC.Range("M62").Value = -2647.4285
$ws_ALC.Range("N62").Value = -3248

$ws_ALC.Range("H65").Value = 3186.6667
$ws_ALC.Range("I65").Value = 3271.4285
$ws_ALC.Range("J65").Value = 2000
$ws_ALC.Range("K65").Value = 16357.1425
$ws_ALC.Range("L65").Value = 10000
$ws_ALC.Range("M65").Value = -13237.1425
$ws_ALC.Range("N65").Value = -16240

$ws_ALC.Range("H100").Value = 2733.6
$ws_ALC.Range("I100").Value = 3334.1667
$ws_ALC.Range("J100").Value = 2333.2222
$ws_ALC.Range("K100").Value = 3334.1667
$ws_ALC.Range("L100").Value = 2333.2222
$ws_ALC.Range("M100").Value = -2793.1667
$ws_ALC.Range("N100").Value = -3415.2222

$ws_ALC.Range("H116").Value = 7812.625
$ws_ALC.Range("I116").Value = 5252.5
$ws_ALC.Range("J116").Value = 8666
$ws_ALC.Range("K116").Value = 5252.5
$ws_ALC.Range("L116").Value = 8666
$ws_ALC.Range("M116").Value = -1810.5
$ws_ALC.Range("N116").Value = -15550

$ws_ALC.Range("H123").Value = 29943
$ws_ALC.Range("J123").Value = 29943
$ws_ALC.Range("L123").Value = 29943
$ws_ALC.Range("N123").Value = -39743

$ws_ALC.Range("H137").Value = 3252.6418
$ws_ALC.Range("I137").Value = 952.4666999999999
$ws_ALC.Range("J137").Value = 3916.1538
$ws_ALC.Range("K137").Value = 2857.4001
$ws_ALC.Range("L137").Value = 11748.4614
$ws_ALC.Range("M137").Value = -307.4000999999998
$ws_ALC.Range("N137").Value = -16848.4614

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H39").Value = 2510.6667
$ws_ARM.Range("I39").Value = 2510.6667
$ws_ARM.Range("K39").Value = 2510.6667
$ws_ARM.Range("M39").Value = -1990.6667

$ws_ARM.Range("H113").Value = 37745.5
$ws_ARM.Range("J113").Value = 37745.5
$ws_ARM.Range("L113").Value = 37745.5
$ws_ARM.Range("N113").Value = -46423.5

$ws_ARM.Range("H122").Value = 2376.5
$ws_ARM.Range("I122").Value = 2544.5715
$ws_ARM.Range("K122").Value = 7633.7145
$ws_ARM.Range("M122").Value = -5183.7145

$ws_ARM.Range("H128").Value = 49996
$ws_ARM.Range("J128").Value = 49996
$ws_ARM.Range("L128").Value = 49996
$ws_ARM.Range("N128").Value = -59956

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H110").Value = 48702
$ws_BSM.Range("J110").Value = 48702
$ws_BSM.Range("L110").Value = 48702
$ws_BSM.Range("N110").Value = -56882

$ws_BSM.Range("H130").Value = 48731.8
$ws_BSM.Range("J130").Value = 48731.8
$ws_BSM.Range("L130").Value = 48731.8
$ws_BSM.Range("N130").Value = -58771.8

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H82").Value = 39566
$ws_CRP.Range("J82").Value = 39566
$ws_CRP.Range("L82").Value = 39566
$ws_CRP.Range("N82").Value = -40288

$ws_CRP.Range("H85").Value = 39566
$ws_CRP.Range("J85").Value = 39566
$ws_CRP.Range("L85").Value = 39566
$ws_CRP.Range("N85").Value = -42062

$ws_CRP.Range("H119").Value = 48757
$ws_CRP.Range("J119").Value = 48757
$ws_CRP.Range("L119").Value = 48757
$ws_CRP.Range("N119").Value = -58433

$ws_CRP.Range("H134").Value = 2247.125
$ws_CRP.Range("I134").Value = 892.875
$ws_CRP.Range("K134").Value = 2678.625
$ws_CRP.Range("M134").Value = -143.625

$ws_CRP.Range("H138").Value = 46451.332
$ws_CRP.Range("J138").Value = 46451.332
$ws_CRP.Range("L138").Value = 46451.332
$ws_CRP.Range("N138").Value = -56731.332

$ws_CRP.Range("H140").Value = 52000
$ws_CRP.Range("J140").Value = 52000
$ws_CRP.Range("L140").Value = 52000
$ws_CRP.Range("N140").Value = -62360

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H122").Value = 2111.1
$ws_GSM.Range("J122").Value = 2334.2222
$ws_GSM.Range("L122").Value = 7002.6666
$ws_GSM.Range("N122").Value = -11902.6666

$ws_GSM.Range("H126").Value = 2482.3635
$ws_GSM.Range("I126").Value = 4103
$ws_GSM.Range("J126").Value = 2122.2222
$ws_GSM.Range("K126").Value = 12309
$ws_GSM.Range("L126").Value = 6366.6666
$ws_GSM.Range("M126").Value = -9839
$ws_GSM.Range("N126").Value = -11306.6666

$ws_GSM.Range("H128").Value = 38441
$ws_GSM.Range("J128").Value = 38441
$ws_GSM.Range("L128").Value = 38441
$ws_GSM.Range("N128").Value = -48401

$ws_GSM.Range("H130").Value = 0
$ws_GSM.Range("J130").Value = 0
$ws_GSM.Range("L130").Value = 0
$ws_GSM.Range("N130").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H36").Value = 46711
$ws_LTW.Range("J36").Value = 46711
$ws_LTW.Range("L36").Value = 46711
$ws_LTW.Range("N36").Value = -47835

$ws_LTW.Range("H93").Value = 3620
$ws_LTW.Range("I93").Value = 3500
$ws_LTW.Range("K93").Value = 3500
$ws_LTW.Range("M93").Value = -2252

$ws_LTW.Range("H137").Value = 24908.334
$ws_LTW.Range("J137").Value = 24908.334
$ws_LTW.Range("L137").Value = 24908.334
$ws_LTW.Range("N137").Value = -35108.334

$ws_LTW.Range("H139").Value = 32963.637
$ws_LTW.Range("I139").Value = 0
$ws_LTW.Range("J139").Value = 32963.637
$ws_LTW.Range("K139").Value = 0
$ws_LTW.Range("L139").Value = 32963.637
$ws_LTW.Range("M139").ClearContents()
$ws_LTW.Range("N139").Value = -43243.637

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H74").Value = 10645.3
$ws_WVR.Range("J74").Value = 10645.3
$ws_WVR.Range("L74").Value = 10645.3
$ws_WVR.Range("N74").Value = -12517.3

$ws_WVR.Range("H75").Value = 30125
$ws_WVR.Range("J75").Value = 30125
$ws_WVR.Range("L75").Value = 30125
$ws_WVR.Range("N75").Value = -31997

$ws_WVR.Range("H77").Value = 10645.3
$ws_WVR.Range("J77").Value = 10645.3
$ws_WVR.Range("L77").Value = 31935.9
$ws_WVR.Range("N77").Value = -41295.89999999999

$ws_WVR.Range("H78").Value = 30125
$ws_WVR.Range("J78").Value = 30125
$ws_WVR.Range("L78").Value = 90375
$ws_WVR.Range("N78").Value = -99735

$ws_WVR.Range("H81").Value = 1854.8
$ws_WVR.Range("I81").Value = 1786.6666
$ws_WVR.Range("J81").Value = 1900.2222
$ws_WVR.Range("K81").Value = 3573.3332
$ws_WVR.Range("L81").Value = 3800.4444
$ws_WVR.Range("M81").Value = -2512.3332
$ws_WVR.Range("N81").Value = -5922.4444

$ws_WVR.Range("H82").Value = 27000
$ws_WVR.Range("J82").Value = 27000
$ws_WVR.Range("L82").Value = 27000
$ws_WVR.Range("N82").Value = -27766

$ws_WVR.Range("H84").Value = 1854.8
$ws_WVR.Range("I84").Value = 1786.6666
$ws_WVR.Range("J84").Value = 1900.2222
$ws_WVR.Range("K84").Value = 17866.666
$ws_WVR.Range("L84").Value = 19002.222
$ws_WVR.Range("M84").Value = -12562.666
$ws_WVR.Range("N84").Value = -29610.222

$ws_WVR.Range("H85").Value = 27000
$ws_WVR.Range("J85").Value = 27000
$ws_WVR.Range("L85").Value = 27000
$ws_WVR.Range("N85").Value = -29652

$ws_WVR.Range("H108").Value = 40311
$ws_WVR.Range("J108").Value = 40311
$ws_WVR.Range("L108").Value = 40311
$ws_WVR.Range("N108").Value = -47991

$ws_WVR.Range("H109").Value = 38369
$ws_WVR.Range("J109").Value = 38369
$ws_WVR.Range("L109").Value = 38369
$ws_WVR.Range("N109").Value = -41143

$ws_WVR.Range("H119").Value = 49694
$ws_WVR.Range("J119").Value = 49694
$ws_WVR.Range("L119").Value = 49694
$ws_WVR.Range("N119").Value = -59370

$ws_WVR.Range("H123").Value = 43473.668
$ws_WVR.Range("J123").Value = 43473.668
$ws_WVR.Range("L123").Value = 43473.668
$ws_WVR.Range("N123").Value = -53273.668

$ws_WVR.Range("H131").Value = 50707
$ws_WVR.Range("J131").Value = 50707
$ws_WVR.Range("L131").Value = 50707
$ws_WVR.Range("N131").Value = -60787

$ws_WVR.Range("H138").Value = 33565.8
$ws_WVR.Range("J138").Value = 33565.8
$ws_WVR.Range("L138").Value = 33565.8
$ws_WVR.Range("N138").Value = -43845.8
